# Added papers for review
#
# Adds a new "Original" data point (a paper result) below the existing
# Nodes/time table on Sheet1, registering a new shared string and a new
# row of data, then leaves the selection where the user ended up after
# typing it in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: "Original" paper result added below the existing "Total" row
$ws.Range("A10").Value = "Original"
$ws.Range("B10").Value = 77.71

# Reflect where the user left the selection after adding the row
$ws.Range("B11").Select() | Out-Null
